$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.531699999999993
$ws.Range("D4").Value = -6.7631
$ws.Range("B7").Value = 4.444200000000002
$ws.Range("A8").Value = -22.25280000000002
$ws.Range("A10").Value = -21.73609999999999
$ws.Range("D11").Value = -8.216499999999995
$ws.Range("A12").Value = -21.5657
$ws.Range("B14").Value = 6.280100000000004
$ws.Range("D14").Value = -7.414199999999998
$ws.Range("B15").Value = 4.380299999999998
$ws.Range("A18").Value = -22.39850000000001
$ws.Range("B18").Value = 4.357399999999996
$ws.Range("D18").Value = -8.16889999999999
$ws.Range("D19").Value = -8.612999999999991
$ws.Range("B20").Value = 9.008499999999996
$ws.Range("D21").Value = -8.274800000000003
$ws.Range("A25").Value = -21.50779999999999
$ws.Range("D27").Value = -8.739700000000001
$ws.Range("B29").Value = 5.063600000000001
$ws.Range("B30").Value = 5.0716
$ws.Range("B31").Value = 5.528300000000004
$ws.Range("D31").Value = -8.236499999999998
$ws.Range("B35").Value = 8.607600000000007
$ws.Range("A37").Value = -20.0794
$ws.Range("D38").Value = -8.7019
$ws.Range("B40").Value = 8.501300000000004
$ws.Range("D42").Value = -9.188599999999994
$ws.Range("B44").Value = 4.611700000000003
$ws.Range("D44").Value = -8.534999999999998
$ws.Range("D47").Value = -7.669499999999999
$ws.Range("B50").Value = 4.338500000000002
$ws.Range("B54").Value = 4.790200000000001
$ws.Range("A55").Value = -21.50840000000001
$ws.Range("D56").Value = -8.403499999999994
$ws.Range("D58").Value = -8.243199999999991
$ws.Range("D65").Value = -7.596599999999999
$ws.Range("A68").Value = -21.48600000000001
$ws.Range("B68").Value = 4.481999999999999
$ws.Range("D73").Value = -7.469199999999998
$ws.Range("B76").Value = 6.343099999999996
$ws.Range("A77").Value = -20.37779999999998
$ws.Range("A78").Value = -20.18539999999997
$ws.Range("A79").Value = -20.03869999999999
$ws.Range("A80").Value = -20.36269999999997
$ws.Range("A81").Value = -21.93289999999999
$ws.Range("A82").Value = -21.6635
$ws.Range("A84").Value = -21.93890000000001
$ws.Range("B87").Value = 4.482599999999994
$ws.Range("B88").Value = 4.324799999999998
$ws.Range("D90").Value = -7.901600000000003
$ws.Range("B92").Value = 4.683799999999997
$ws.Range("D92").Value = -6.189199999999999
$ws.Range("D94").Value = -6.794499999999999
$ws.Range("D95").Value = -8.033599999999998
$ws.Range("B96").Value = 4.678400000000005
$ws.Range("B98").Value = 6.0678
$ws.Range("A101").Value = -21.69
$ws.Range("B101").Value = 5.735500000000004
$ws.Range("D101").Value = -7.834899999999997
$ws.Range("A102").Value = -21.97429999999999
$ws.Range("B102").Value = 4.782400000000002
